$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"3.8487633065266E-19"
$ws.Range("E2").Value = [double]"3.8487633065266E-19"

$ws.Range("D3").Value = [double]"0.9999999635494419"
$ws.Range("E3").Value = [double]"0.9999999635494419"

$ws.Range("D4").Value = [double]"0.9999999999213351"
$ws.Range("E4").Value = [double]"7.866485240981547E-11"

$ws.Range("D5").Value = [double]"0.635671400665928"
$ws.Range("E5").Value = [double]"0.364328599334072"

$ws.Range("D6").Value = [double]"0.9986787667894935"
$ws.Range("E6").Value = [double]"0.001321233210506478"

$ws.Range("D8").Value = [double]"0.9963973763030217"
$ws.Range("E8").Value = [double]"0.00360262369697828"
$ws.Range("F8").Value = [double]"2.512187957763672"
